# Updated cryptos list on Wed May  3 03:29:23 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.618.07"

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.870.11"
$ws.Range("E3").Value = "  +2.26%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.30%  "

# Row 5 - BNB
$ws.Range("D5").Value = "326.58"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.27%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4655"

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3884"
$ws.Range("E8").Value = "  +0.64%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.07870"
$ws.Range("E9").Value = "  +0.18%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "0.9726"
$ws.Range("E10").Value = "  +1.67%  "

# Row 11 - Solana
$ws.Range("D11").Value = "21.96"
$ws.Range("E11").Value = "  +0.67%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.819.58"
$ws.Range("E12").Value = "  -4.81%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "6.988"
$ws.Range("E13").Value = "  +1.62%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "5.696"
$ws.Range("E14").Value = "  +0.85%  "

# Row 15 - TRON
$ws.Range("D15").Value = "0.06977"
$ws.Range("E15").Value = "  +3.36%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "88.06"
$ws.Range("E16").Value = "  +1.50%  "

# Row 17 - BinanceUSD
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  +0.44%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.00001005"
$ws.Range("E18").Value = "  +1.42%  "

# Row 19 - Avalanche
$ws.Range("D19").Value = "16.80"
$ws.Range("E19").Value = "  +1.39%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.34%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "28.618.33"
$ws.Range("E21").Value = "  +2.12%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "5.291"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "10.99"
$ws.Range("E23").Value = "  +0.33%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.116"

# Row 25 - WrappedliquidstakedEther2.0
$ws.Range("D25").Value = "2.108.28"
$ws.Range("E25").Value = "  -0.14%  "

# Row 26 - Monero
$ws.Range("D26").Value = "152.79"
$ws.Range("E26").Value = "  -0.66%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "19.21"
$ws.Range("E27").Value = "  +0.57%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "5.769"
$ws.Range("E28").Value = "  +0.92%  "

# Row 29 - LidoDAOToken
$ws.Range("E29").Value = "  +0.92%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "119.32"
$ws.Range("E30").Value = "  +2.03%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.09371"
$ws.Range("E31").Value = "  +1.45%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "0.9195"
$ws.Range("E32").Value = "  -1.61%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "5.270"
$ws.Range("E33").Value = "  -0.28%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  +1.92%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "3.340"
$ws.Range("E35").Value = "  +0.63%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "0.05800"
$ws.Range("E36").Value = "  -0.92%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -2.16%  "

# Row 38 - TrustWalletToken
$ws.Range("E38").Value = "  +0.16%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "7.744"
$ws.Range("E39").Value = "  -0.16%  "

# Row 40 - TheSandbox
$ws.Range("E40").Value = "  +0.92%  "

# Row 41 - Algorand
$ws.Range("E41").Value = "  +1.52%  "

# Row 42 - Aptos
$ws.Range("E42").Value = "  -1.33%  "

# Row 43 - Cronos
$ws.Range("D43").Value = "0.07210"
$ws.Range("E43").Value = "  +2.74%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "11.71"
$ws.Range("E44").Value = "  +1.38%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "0.5309"
$ws.Range("E45").Value = "  +1.08%  "

# Row 46 - WEMIXToken
$ws.Range("E46").Value = "  -6.34%  "

# Row 47 - now RenderToken (was Quant)
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "2.080"
$ws.Range("E47").Value = "  -2.65%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "1.821"
$ws.Range("E48").Value = "  -0.21%  "

# Row 49 - now Quant (was RenderToken)
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "113.11"
$ws.Range("E49").Value = "  +0.50%  "

# Row 50 - MXToken
$ws.Range("D50").Value = "2.407"
$ws.Range("E50").Value = "  +3.79%  "

# Row 51 - PaxDollar
$ws.Range("D51").Value = "1.005"
$ws.Range("E51").Value = "  +0.39%  "
